$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 897094.9399999999
$ws.Range("I17").Value = 363
$ws.Range("J17").Value = 909212.9399999999
$ws.Range("K17").Value = 1089
$ws.Range("L17").Value = 2727638.82
$ws.Range("M17").Value = -921
$ws.Range("N17").Value = -2727974.82
$ws.Range("H75").Value = 118328.25
$ws.Range("J75").Value = 118328.25
$ws.Range("L75").Value = 118328.25
$ws.Range("N75").Value = -120200.25
$ws.Range("H76").Value = 20806.285
$ws.Range("J76").Value = 6969
$ws.Range("L76").Value = 6969
$ws.Range("N76").Value = -7599
$ws.Range("H78").Value = 118328.25
$ws.Range("J78").Value = 118328.25
$ws.Range("L78").Value = 354984.75
$ws.Range("N78").Value = -364344.75
$ws.Range("H79").Value = 20806.285
$ws.Range("J79").Value = 6969
$ws.Range("L79").Value = 6969
$ws.Range("N79").Value = -9153
$ws.Range("H131").Value = 3541.889
$ws.Range("I131").Value = 3109.875
$ws.Range("J131").Value = 6998
$ws.Range("K131").Value = 9329.625
$ws.Range("L131").Value = 20994
$ws.Range("M131").Value = -4289.625
$ws.Range("N131").Value = -31074
$ws.Range("H138").Value = 14927972
$ws.Range("I138").Value = 1000.9231
$ws.Range("K138").Value = 3002.7693
$ws.Range("M138").Value = 2137.2307

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 58139.57
$ws.Range("I74").Value = 60721.55
$ws.Range("J74").Value = 6500
$ws.Range("K74").Value = 60721.55
$ws.Range("L74").Value = 6500
$ws.Range("M74").Value = -59847.55
$ws.Range("N74").Value = -8248
$ws.Range("H77").Value = 58139.57
$ws.Range("I77").Value = 60721.55
$ws.Range("J77").Value = 6500
$ws.Range("K77").Value = 303607.75
$ws.Range("L77").Value = 32500
$ws.Range("M77").Value = -299239.75
$ws.Range("N77").Value = -41236
$ws.Range("H132").Value = 75072.32000000001
$ws.Range("I132").Value = 4843.793
$ws.Range("K132").Value = 14531.379
$ws.Range("M132").Value = -12001.379

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1551.75
$ws.Range("I20").Value = 1363.8889
$ws.Range("J20").Value = 1889.9
$ws.Range("K20").Value = 1363.8889
$ws.Range("L20").Value = 1889.9
$ws.Range("M20").Value = -1116.8889
$ws.Range("N20").Value = -2383.9
$ws.Range("H92").Value = 276741
$ws.Range("J92").Value = 276741
$ws.Range("L92").Value = 276741
$ws.Range("N92").Value = -281733
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2981.611
$ws.Range("I58").Value = 4133.6665
$ws.Range("J58").Value = 1829.5555
$ws.Range("K58").Value = 4133.6665
$ws.Range("L58").Value = 1829.5555
$ws.Range("M58").Value = -3930.6665
$ws.Range("N58").Value = -2235.5555
$ws.Range("H131").Value = 65100
$ws.Range("J131").Value = 65100
$ws.Range("L131").Value = 65100
$ws.Range("N131").Value = -75180
$ws.Range("H134").Value = 2461.0908
$ws.Range("I134").Value = 2298.2856
$ws.Range("K134").Value = 6894.8568
$ws.Range("M134").Value = -4359.8568
$ws.Range("H136").Value = 2981.611
$ws.Range("I136").Value = 4133.6665
$ws.Range("J136").Value = 1829.5555
$ws.Range("K136").Value = 12400.9995
$ws.Range("L136").Value = 5488.666499999999
$ws.Range("M136").Value = -9850.999500000002
$ws.Range("N136").Value = -10588.6665
$ws.Range("H141").Value = 86802.19
$ws.Range("J141").Value = 86802.19
$ws.Range("L141").Value = 86802.19
$ws.Range("N141").Value = -97162.19

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H99").Value = 3446.5715
$ws.Range("I99").Value = 3825
$ws.Range("K99").Value = 11475
$ws.Range("M99").Value = -9229
$ws.Range("H128").Value = 159900
$ws.Range("I128").Value = 159900
$ws.Range("K128").Value = 479700
$ws.Range("M128").Value = -474720
$ws.Range("H131").Value = 1630.35
$ws.Range("J131").Value = 1743.0714
$ws.Range("L131").Value = 5229.2142
$ws.Range("N131").Value = -15309.2142

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 50000
$ws.Range("I42").Value = 50000
$ws.Range("K42").Value = 50000
$ws.Range("M42").Value = -49515
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 830.9167
$ws.Range("I97").Value = 638.9231
$ws.Range("J97").Value = 1057.8182
$ws.Range("K97").Value = 638.9231
$ws.Range("L97").Value = 1057.8182
$ws.Range("M97").Value = -142.9231
$ws.Range("N97").Value = -2049.8182
$ws.Range("H102").Value = 4204.724
$ws.Range("I102").Value = 2287.158
$ws.Range("K102").Value = 2287.158
$ws.Range("M102").Value = -665.1579999999999
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H112").Value = 72000
$ws.Range("J112").Value = 72000
$ws.Range("L112").Value = 72000
$ws.Range("N112").Value = -74216
$ws.Range("H113").Value = 2991.8572
$ws.Range("I113").Value = 1943.9
$ws.Range("K113").Value = 1943.9
$ws.Range("M113").Value = 226.0999999999999
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H115").Value = 50000
$ws.Range("I115").Value = 50000
$ws.Range("K115").Value = 50000
$ws.Range("M115").Value = -48825

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2686.5908
$ws.Range("I82").Value = 2593.7144
$ws.Range("K82").Value = 2593.7144
$ws.Range("M82").Value = -2232.7144
$ws.Range("H85").Value = 2686.5908
$ws.Range("I85").Value = 2593.7144
$ws.Range("K85").Value = 2593.7144
$ws.Range("M85").Value = -1345.7144
$ws.Range("H122").Value = 5366
$ws.Range("I122").Value = 4452.6
$ws.Range("J122").Value = 7649.5
$ws.Range("K122").Value = 13357.8
$ws.Range("L122").Value = 22948.5
$ws.Range("M122").Value = -10907.8
$ws.Range("N122").Value = -27848.5
$ws.Range("H131").Value = 89999
$ws.Range("J131").Value = 89999
$ws.Range("L131").Value = 89999
$ws.Range("N131").Value = -100079
$ws.Range("H132").Value = 2004.0294
$ws.Range("I132").Value = 1811.1936
$ws.Range("K132").Value = 5433.5808
$ws.Range("M132").Value = -2903.5808

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7249.5
$ws.Range("I62").Value = 5499
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 5499
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -4875
$ws.Range("N62").Value = -10248
$ws.Range("H65").Value = 7249.5
$ws.Range("I65").Value = 5499
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 27495
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -24375
$ws.Range("N65").Value = -51240
$ws.Range("H100").Value = 30304130
$ws.Range("I100").Value = 43479400
$ws.Range("K100").Value = 86958800
$ws.Range("M100").Value = -86958259
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 2936.1052
$ws.Range("I136").Value = 2013.7567
$ws.Range("J136").Value = 4642.45
$ws.Range("K136").Value = 6041.2701
$ws.Range("L136").Value = 13927.35
$ws.Range("M136").Value = -3491.2701
$ws.Range("N136").Value = -19027.35
